function Find-ParaIndex($doc, $prefix) {
    $i = 0
    foreach ($para in $doc.Paragraphs) {
        $i = $i + 1
        if ($para.Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Test cases:" heading paragraph gets a bold paragraph mark,
#    and the following empty paragraph becomes two new paragraphs:
#      "Prerequisites: " (with "Prerequisites:" bold)
#      " User is logged in."
# ---------------------------------------------------------------
$idxTestCases = Find-ParaIndex $d "Test cases:"
$pTestCases = $d.Paragraphs($idxTestCases)
$rTestCases = $pTestCases.Range

# Bold the whole paragraph (text + trailing paragraph mark), then remove
# the bold from the visible text only, leaving the paragraph-mark bold.
$rTestCases.Font.Bold = 1
$textOnly = $d.Range($rTestCases.Start, $rTestCases.End - 1)
$textOnly.Font.Bold = 0

# The paragraph right after "Test cases:" is the empty placeholder paragraph.
$idxEmpty = $idxTestCases + 1
$pEmpty = $d.Paragraphs($idxEmpty)
$rEmpty = $pEmpty.Range
$rEmpty.InsertBefore("Prerequisites: ")

# Split this paragraph into its own paragraph, then add the next sentence.
$pPrereq = $d.Paragraphs($idxEmpty)
$rPrereq = $pPrereq.Range
$rPrereq.InsertParagraphAfter()

$idxUserLoggedIn = $idxEmpty + 1
$pUserLoggedIn = $d.Paragraphs($idxUserLoggedIn)
$rUserLoggedIn = $pUserLoggedIn.Range
$rUserLoggedIn.InsertBefore(" User is logged in.")

# Bold just the "Prerequisites:" word (not the trailing space) in its paragraph.
$pPrereq2 = $d.Paragraphs($idxEmpty)
$rPrereq2 = $pPrereq2.Range
$boldRange = $d.Range($rPrereq2.Start, $rPrereq2.Start + "Prerequisites:".Length)
$boldRange.Font.Bold = 1
$boldRange.Font.BoldBi = 1

# ---------------------------------------------------------------
# 2) Merge the two runs that make up
#    "User clicks on confirm, along with submitting information failed"
#    into a single run, without touching the separate "2.2: " run.
# ---------------------------------------------------------------
$idx22 = Find-ParaIndex $d "2.2: "
$p22 = $d.Paragraphs($idx22)
$r22 = $p22.Range

$prefixLen = "2.2: ".Length
$run1 = $d.Range($r22.Start, $r22.Start + $prefixLen)

# Temporarily bold the "2.2: " run so it keeps its own identity while the
# rest of the paragraph's runs get rewritten/merged below.
$run1.Font.Bold = 1

$p22b = $d.Paragraphs($idx22)
$r22b = $p22b.Range
$mergedText = "User clicks on confirm, along with submitting information failed"
$sub = $d.Range($r22b.Start + $prefixLen, $r22b.End - 1)
$sub.Delete()

$p22c = $d.Paragraphs($idx22)
$r22c = $p22c.Range
$insPoint = $d.Range($r22c.End - 1, $r22c.End - 1)
$insPoint.InsertBefore($mergedText)

# Clear the bold on the newly (re)created merged run.
$p22d = $d.Paragraphs($idx22)
$r22d = $p22d.Range
$newRun = $d.Range($r22d.Start + $prefixLen, $r22d.End - 1)
$newRun.Font.Bold = 0

# Clear the temporary bold flag on "2.2: " again.
$run1b = $d.Range($r22.Start, $r22.Start + $prefixLen)
$run1b.Font.Bold = 0
